# Swap the data values between row 2 and row 3 for the columns that
# actually differ (A, B, D, E, F, G, H, Q, R). Columns that are identical
# between the two rows (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AN, AO, AT, AW, AX, AY) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"

    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2

    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
